$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 70

$ws.Cells.Item($row, 1).Value = "2025-08-29 03:46:58 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-29 09:16:58 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""

$rng = $ws.Range("A70:H70")
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
